$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data row: userId 4 -> 5, name "Zaid" -> "Zaidi"
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Zaidi"

# Move the active selection to B2 (as recorded in the saved view state)
$ws.Range("B2").Select()
